# ------------------------------------------------------------------
# Applies the edits described by the diff:
#  1. Clears the (already-empty) B7 / B13 cells on "ODI Batting" so the
#     <c> element is dropped entirely (matches the diff removing them).
#  2. Adds a new worksheet "ODI Batting Extra" (4th sheet, after
#     "ODI Bowling") containing the MATCH_CODE / BATTING_POSITION /
#     NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL / MAN_OF_MATCH table.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. "ODI Batting" sheet: drop the two stray empty B7 / B13 cells ---
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("B7").ClearContents()
$odiBatting.Range("B13").ClearContents()

# --- 2. Add the new "ODI Batting Extra" worksheet as the last tab ------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$extra.Name = "ODI Batting Extra"

# Match the outline / page-margin conventions used by the other sheets
# in this workbook.
$extra.Outline.SummaryRow = 1
$extra.Outline.SummaryColumn = 1
$extra.PageSetup.LeftMargin = 54
$extra.PageSetup.RightMargin = 54
$extra.PageSetup.TopMargin = 72
$extra.PageSetup.BottomMargin = 72
$extra.PageSetup.HeaderMargin = 36
$extra.PageSetup.FooterMargin = 36

# Header row (copy the bold/centered/bordered header style used on the
# other sheets so the new header matches visually).
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($col = 1; $col -le $headers.Length; $col++) {
    $extra.Cells.Item(1, $col).Value = $headers[$col - 1]
}
$odiBatting.Range("A1:F1").Copy()
$extra.Range("A1:F1").PasteSpecial(-4122)   # xlPasteFormats

# Helper: write a value into a cell as TEXT (so numeric-looking strings
# like "4457" or "11.18%" are stored as strings, not coerced to numbers),
# without leaving any custom number-format style attached to the cell.
function Set-TextCell($cell, $text) {
    if ($text -eq $null -or $text -eq "") {
        return
    }
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Data rows (row 2 .. row 13). $null marks a blank cell.
$rows = @(
    @("4456", $null, $null, $null, $null,    "NO"),
    @("4457", 6,     "4",   "1",   "11.18%", "NO"),
    @("4469", 2,     "0",   "1",   "4.76%",  "NO"),
    @("4598", $null, $null, $null, $null,    "NO"),
    @("4599", $null, $null, $null, $null,    "NO"),
    @("4602", 5,     $null, $null, $null,    "NO"),
    @("4609", 6,     "0",   "0",   $null,    "NO"),
    @("4613", 6,     "2",   "2",   "13.41%", "NO"),
    @("4618", 7,     "2",   "2",   "10.42%", "NO"),
    @("4619", $null, $null, $null, $null,    "NO"),
    @("4620", 7,     "1",   "3",   "18.91%", "NO"),
    @("4622", $null, $null, $null, $null,    "NO")
)

$r = 2
foreach ($row in $rows) {
    Set-TextCell $extra.Cells.Item($r, 1) $row[0]               # MATCH_CODE (text)

    if ($row[1] -ne $null) {
        $extra.Cells.Item($r, 2).Value = $row[1]                # BATTING_POSITION (real number)
    }

    Set-TextCell $extra.Cells.Item($r, 3) $row[2]                # NUM_4 (text)
    Set-TextCell $extra.Cells.Item($r, 4) $row[3]                # NUM_6 (text)
    Set-TextCell $extra.Cells.Item($r, 5) $row[4]                # PERCENT_RUNS_OF_TOTAL (text)
    Set-TextCell $extra.Cells.Item($r, 6) $row[5]                # MAN_OF_MATCH (text)

    $r = $r + 1
}

# Restore the original active sheet/selection ("Player Info" was active
# before this edit) now that the new sheet has been populated.
$wb.Worksheets.Item("Player Info").Activate() | Out-Null
$wb.Worksheets.Item("Player Info").Range("A1").Select() | Out-Null
